# Clean-up of input table: the workbook was re-saved from a different
# Excel install (Windows, German locale defaults) which renamed the
# worksheet tab back to its default name and left the cursor on a
# different cell. Reproduce the two user-visible, content-level changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "updated" to "Tabelle1"
$ws.Name = "Tabelle1"

# Update the active selection to C9 (was C12)
$ws.Range("C9").Select()
